$d = $word.ActiveDocument

# The document has a series of "Normal (Web)" paragraphs whose paragraph
# shading (w:shd) currently fills with a light grey (#F4F5F7). The edit
# changes that shading to white (#FFFFFF), matching Word's
# "White, Background 1" theme swatch.
$oldFill = 16250356   # BGR-encoded RGB(0xF4,0xF5,0xF7)
$newFill = 16777215   # RGB(255,255,255) / wdColorWhite

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $shading = $para.Range.Shading
    if ($shading.BackgroundPatternColor -eq $oldFill) {
        $shading.BackgroundPatternColor = $newFill
    }
}

Write-Output "shading updated"
